# Apply January 2024 data-refresh figures to the department / agent-handler
# summary sheets. Values in this workbook are stored as plain text (they were
# authored by an external export tool as inline strings, not real numbers),
# so every write below is prefixed with a leading apostrophe to force Excel
# to keep them as text instead of silently re-typing them as numeric values
# (which would also reformat things like "3000.00" down to "3000").

$wb = $excel.ActiveWorkbook

# --- Sheet "部门情况202401" (department overview) --------------------------
$wsDept = $wb.Sheets.Item("部门情况202401")

$wsDept.Range("J2").Value = "'152624.45"
$wsDept.Range("J3").Value = "'11145.17"
$wsDept.Range("J4").Value = "'63578.19"
$wsDept.Range("J5").Value = "'3000.00"
$wsDept.Range("J6").Value = "'0.00"
$wsDept.Range("J7").Value = "'89381.07"

# --- Sheet "经办人情况202401" (handler / agent overview) --------------------
$wsAgent = $wb.Sheets.Item("经办人情况202401")

$wsAgent.Range("E3").Value = "'0.00"
$wsAgent.Range("F3").Value = "'0.00"

$wsAgent.Range("E5").Value = "'0.00"
$wsAgent.Range("F5").Value = "'0.00"

$wsAgent.Range("E6").Value = "'0.00"
$wsAgent.Range("F6").Value = "'0.00"

$wsAgent.Range("E7").Value = "'0.00"
$wsAgent.Range("F7").Value = "'0.00"

$wsAgent.Range("E8").Value = "'0.00"
$wsAgent.Range("F8").Value = "'0.00"

$wsAgent.Range("E9").Value = "'37.61"
$wsAgent.Range("F9").Value = "'2.00"

$wsAgent.Range("E10").Value = "'11107.56"
$wsAgent.Range("F10").Value = "'22.00"

$wsAgent.Range("E11").Value = "'0.00"
$wsAgent.Range("F11").Value = "'0.00"

$wsAgent.Range("E12").Value = "'6411.81"
$wsAgent.Range("F12").Value = "'22.00"

$wsAgent.Range("E13").Value = "'97.82"
$wsAgent.Range("F13").Value = "'1.00"

$wsAgent.Range("E14").Value = "'29559.88"
$wsAgent.Range("F14").Value = "'34.00"

$wsAgent.Range("E15").Value = "'0.00"
$wsAgent.Range("F15").Value = "'0.00"

$wsAgent.Range("E16").Value = "'0.00"
$wsAgent.Range("F16").Value = "'0.00"

$wsAgent.Range("E17").Value = "'3000.00"
$wsAgent.Range("F17").Value = "'6.00"

$wsAgent.Range("E18").Value = "'0.00"
$wsAgent.Range("F18").Value = "'0.00"

$wsAgent.Range("E19").Value = "'0.00"
$wsAgent.Range("F19").Value = "'0.00"

$wsAgent.Range("E27").Value = "'0.00"
$wsAgent.Range("F27").Value = "'0.00"
